$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 1497.5
$ws.Range("I74").Value = 1497.5
$ws.Range("K74").Value = 1497.5
$ws.Range("M74").Value = -561.5
$ws.Range("H77").Value = 1497.5
$ws.Range("I77").Value = 1497.5
$ws.Range("K77").Value = 7487.5
$ws.Range("M77").Value = -2807.5
$ws.Range("H88").Value = 13475.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 13475.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 13475.5
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -14287.5
$ws.Range("H91").Value = 13475.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 13475.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 13475.5
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -16283.5
$ws.Range("H98").Value = 3379.1333
$ws.Range("I98").Value = 3275.75
$ws.Range("K98").Value = 3275.75
$ws.Range("M98").Value = -1777.75
$ws.Range("H122").Value = 3379.1333
$ws.Range("I122").Value = 3275.75
$ws.Range("K122").Value = 9827.25
$ws.Range("M122").Value = -7377.25
$ws.Range("H132").Value = 2885.9092
$ws.Range("I132").Value = 2826.5789
$ws.Range("J132").Value = 3261.6667
$ws.Range("K132").Value = 8479.736699999999
$ws.Range("L132").Value = 9785.000100000001
$ws.Range("M132").Value = -5949.736699999999
$ws.Range("N132").Value = -14845.0001
$ws.Range("H137").Value = 2145.2666
$ws.Range("I137").Value = 1698.4166
$ws.Range("K137").Value = 5095.2498
$ws.Range("M137").Value = -2545.2498
$ws.Range("H138").Value = 3015.7896
$ws.Range("I138").Value = 1491.75
$ws.Range("K138").Value = 4475.25
$ws.Range("M138").Value = 664.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2678
$ws.Range("I45").Value = 2443.8572
$ws.Range("J45").Value = 3497.5
$ws.Range("K45").Value = 2443.8572
$ws.Range("L45").Value = 3497.5
$ws.Range("M45").Value = -2066.8572
$ws.Range("N45").Value = -4251.5
$ws.Range("H61").Value = 1760.5714
$ws.Range("I61").Value = 1760.5714
$ws.Range("K61").Value = 1760.5714
$ws.Range("M61").Value = -1548.5714
$ws.Range("H74").Value = 1354.8
$ws.Range("I74").Value = 1132.5385
$ws.Range("K74").Value = 1132.5385
$ws.Range("M74").Value = -258.5385000000001
$ws.Range("H77").Value = 1354.8
$ws.Range("I77").Value = 1132.5385
$ws.Range("K77").Value = 5662.692500000001
$ws.Range("M77").Value = -1294.692500000001
$ws.Range("H110").Value = 3260.35
$ws.Range("I110").Value = 1912
$ws.Range("J110").Value = 4363.5454
$ws.Range("K110").Value = 1912
$ws.Range("L110").Value = 4363.5454
$ws.Range("M110").Value = 133
$ws.Range("N110").Value = -8453.545399999999
$ws.Range("H132").Value = 3593.5
$ws.Range("I132").Value = 3203.6843
$ws.Range("K132").Value = 9611.052899999999
$ws.Range("M132").Value = -7081.052899999999
$ws.Range("H136").Value = 1760.5714
$ws.Range("I136").Value = 1760.5714
$ws.Range("K136").Value = 5281.7142
$ws.Range("M136").Value = -2731.7142

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3010.7778
$ws.Range("I20").Value = 2038.8
$ws.Range("J20").Value = 4225.75
$ws.Range("K20").Value = 2038.8
$ws.Range("L20").Value = 4225.75
$ws.Range("M20").Value = -1791.8
$ws.Range("N20").Value = -4719.75
$ws.Range("H94").Value = 2272.2104
$ws.Range("I94").Value = 2287.611
$ws.Range("J94").Value = 1995
$ws.Range("K94").Value = 2287.611
$ws.Range("L94").Value = 1995
$ws.Range("M94").Value = -1836.611
$ws.Range("N94").Value = -2897
$ws.Range("H134").Value = 12324.5
$ws.Range("I134").Value = 14299.333
$ws.Range("J134").Value = 6400
$ws.Range("K134").Value = 42897.999
$ws.Range("L134").Value = 19200
$ws.Range("M134").Value = -40362.999
$ws.Range("N134").Value = -24270

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7550
$ws.Range("I58").Value = 11366.667
$ws.Range("J58").Value = 1825
$ws.Range("K58").Value = 11366.667
$ws.Range("L58").Value = 1825
$ws.Range("M58").Value = -11163.667
$ws.Range("N58").Value = -2231
$ws.Range("H63").Value = 70000
$ws.Range("J63").Value = 70000
$ws.Range("L63").Value = 70000
$ws.Range("N63").Value = -71372
$ws.Range("H66").Value = 70000
$ws.Range("J66").Value = 70000
$ws.Range("L66").Value = 210000
$ws.Range("N66").Value = -216864
$ws.Range("H99").Value = 2816.5
$ws.Range("I99").Value = 1779.8
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 1779.8
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = -281.8
$ws.Range("N99").Value = -10996
$ws.Range("H107").Value = 634.125
$ws.Range("I107").Value = 483
$ws.Range("J107").Value = 724.8
$ws.Range("K107").Value = 483
$ws.Range("L107").Value = 724.8
$ws.Range("M107").Value = 1437
$ws.Range("N107").Value = -4564.8
$ws.Range("H109").Value = 56996.75
$ws.Range("J109").Value = 56996.75
$ws.Range("L109").Value = 56996.75
$ws.Range("N109").Value = -59076.75
$ws.Range("H126").Value = 2816.5
$ws.Range("I126").Value = 1779.8
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 5339.4
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -2869.4
$ws.Range("N126").Value = -28940
$ws.Range("H134").Value = 1833.4445
$ws.Range("I134").Value = 2035.5714
$ws.Range("K134").Value = 6106.7142
$ws.Range("M134").Value = -3571.7142
$ws.Range("H136").Value = 7550
$ws.Range("I136").Value = 11366.667
$ws.Range("J136").Value = 1825
$ws.Range("K136").Value = 34100.001
$ws.Range("L136").Value = 5475
$ws.Range("M136").Value = -31550.001
$ws.Range("N136").Value = -10575

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1405.5454
$ws.Range("I4").Value = 913.4761999999999
$ws.Range("J4").Value = 2266.6667
$ws.Range("K4").Value = 2740.4286
$ws.Range("L4").Value = 6800.000100000001
$ws.Range("M4").Value = -2628.4286
$ws.Range("N4").Value = -7024.000100000001
$ws.Range("H23").Value = 557.2
$ws.Range("J23").Value = 452.44446
$ws.Range("L23").Value = 1357.33338
$ws.Range("N23").Value = -1827.33338
$ws.Range("H33").Value = 98.25
$ws.Range("H104").Value = 19833.334
$ws.Range("I104").Value = 14000
$ws.Range("K104").Value = 42000
$ws.Range("M104").Value = -39379
$ws.Range("H132").Value = 1475.8334
$ws.Range("I132").Value = 1071
$ws.Range("K132").Value = 9639
$ws.Range("M132").Value = -7109

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H99").Value = 10982.625
$ws.Range("I99").Value = 12532.2
$ws.Range("J99").Value = 8400
$ws.Range("K99").Value = 12532.2
$ws.Range("L99").Value = 8400
$ws.Range("M99").Value = -10286.2
$ws.Range("N99").Value = -12892
$ws.Range("H102").Value = 2042.4546
$ws.Range("I102").Value = 2042.4546
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2042.4546
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -420.4546
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 1163.6666
$ws.Range("I132").Value = 996.2857
$ws.Range("J132").Value = 1749.5
$ws.Range("K132").Value = 2988.8571
$ws.Range("L132").Value = 5248.5
$ws.Range("M132").Value = -458.8571000000002
$ws.Range("N132").Value = -10308.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 8000
$ws.Range("J38").Value = 8000
$ws.Range("L38").Value = 8000
$ws.Range("N38").Value = -8820
$ws.Range("H55").Value = 1195.7142
$ws.Range("I55").Value = 1274
$ws.Range("K55").Value = 1274
$ws.Range("M55").Value = -1101
$ws.Range("H68").Value = 4500
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 4500
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51352

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 876.2857
$ws.Range("I100").Value = 827.4
$ws.Range("K100").Value = 1654.8
$ws.Range("M100").Value = -1113.8
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
